$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q4" worksheet right after "总计" by duplicating
#    the existing "2022-Q3" sheet (same column layout/format), then
#    overwrite its data with the 2022-Q4 numbers. All other quarter
#    tabs ("2022-Q3","2022-Q2","2022-Q1","2021-Q4") simply shift one
#    slot to the right - no edits needed on them.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q3src = $wb.Worksheets.Item(2)
$q3src.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Force the text-typed columns (B:G) to stay text so numeric-looking
# strings (fund code, pct figures) are not reinterpreted as numbers.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "016105"
$q4.Range("C2").Value = "申万菱信兴乐优选混合A"
$q4.Range("D2").Value = "3.37"
$q4.Range("E2").Value = "54.93"
$q4.Range("F2").Value = "2.70"
$q4.Range("G2").Value = "0.0910"
$q4.Range("H2").Value = 8

$q4.Range("B3").Value = "016106"
$q4.Range("C3").Value = "申万菱信兴乐优选混合C"
$q4.Range("D3").Value = "2.90"
$q4.Range("E3").Value = "54.93"
$q4.Range("F3").Value = "2.70"
$q4.Range("G3").Value = "0.0783"
$q4.Range("H3").Value = 8

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a 2022-Q4 row at the
#    top of the data (row 2) and push the existing quarters down one
#    row each (row2->3, row3->4, row4->5, row5->6).
# ------------------------------------------------------------------

# Row 5 -> Row 6 (2021-Q4 data). Copy A5's format into A6 first so the
# new index cell matches the rest of column A, then copy the row values.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)
$total.Range("B5:D5").Copy()
$total.Range("B6:D6").PasteSpecial(-4163)
$total.Range("A6").Value = 4

# Row 4 -> Row 5 (2022-Q1 data)
$total.Range("B4:D4").Copy()
$total.Range("B5:D5").PasteSpecial(-4163)

# Row 3 -> Row 4 (2022-Q2 data)
$total.Range("B3:D3").Copy()
$total.Range("B4:D4").PasteSpecial(-4163)

# Row 2 -> Row 3 (2022-Q3 data)
$total.Range("B2:D2").Copy()
$total.Range("B3:D3").PasteSpecial(-4163)

# New Row 2 (2022-Q4 data)
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

# ------------------------------------------------------------------
# 3) Restore "总计" as the active sheet (it was active before the edit).
# ------------------------------------------------------------------
$total.Activate()
